$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.322.06'
$ws.Range("E2").Value = '  -0.93%  '
$ws.Range("D3").Value = '1.872.41'
$ws.Range("E3").Value = '  -0.24%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.78'
$ws.Range("E5").Value = '  -1.64%  '
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4692'
$ws.Range("E7").Value = '  -1.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2875'
$ws.Range("E8").Value = '  -0.95%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06444'
$ws.Range("E9").Value = '  -0.59%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.03'
$ws.Range("E10").Value = '  +0.15%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07768'
$ws.Range("E11").Value = '  +0.46%  '
$ws.Range("D12").Value = '1.876.58'
$ws.Range("E12").Value = '  -0.04%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '95.65'
$ws.Range("E13").Value = '  -0.29%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7208'
$ws.Range("E14").Value = '  -2.57%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.123'
$ws.Range("E15").Value = '  -1.15%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '279.34'
$ws.Range("E16").Value = '  +1.72%  '
$ws.Range("D17").Value = '30.310.35'
$ws.Range("E17").Value = '  -0.94%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.98'
$ws.Range("E18").Value = '  -2.10%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.000'
$ws.Range("E19").Value = '  +0.02%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007425'
$ws.Range("E20").Value = '  -0.69%  '
$ws.Range("D21").Value = '2.122.74'
$ws.Range("E21").Value = '  +0.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9997'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.219'
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.225'
$ws.Range("E24").Value = '  +0.84%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '163.23'
$ws.Range("E25").Value = '  -1.39%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.042'
$ws.Range("E26").Value = '  -1.63%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.66'
$ws.Range("E27").Value = '  -0.46%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.879'
$ws.Range("E28").Value = '  -1.23%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.320'
$ws.Range("E29").Value = '  -1.86%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09560'
$ws.Range("E30").Value = '  -3.36%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.469'
$ws.Range("E31").Value = '  -2.65%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.217'
$ws.Range("E33").Value = '  +0.19%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04805'
$ws.Range("E34").Value = '  +0.56%  '
$ws.Range("E35").Value = '  -0.16%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6865'
$ws.Range("E36").Value = '  -1.06%  '
$ws.Range("E38").Value = '  +0.80%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.810'
$ws.Range("E39").Value = '  +1.79%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.218'
$ws.Range("E40").Value = '  -0.14%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '74.16'
$ws.Range("E41").Value = '  +1.11%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4232'
$ws.Range("E42").Value = '  +1.71%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.928'
$ws.Range("E43").Value = '  -2.06%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8253'
$ws.Range("E45").Value = '  -1.12%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '100.65'
$ws.Range("E46").Value = '  -0.89%  '
$ws.Range("E47").Value = '  +2.38%  '
$ws.Range("B48").Value = 'Elrond'
$ws.Range("C48").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '35.07'
$ws.Range("E48").Value = '  -0.66%  '
$ws.Range("B49").Value = 'Aptos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.915'
$ws.Range("E49").Value = '  -0.72%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '897.01'
$ws.Range("E50").Value = '  -2.43%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05717'
$ws.Range("E51").Value = '  +0.84%  '
